$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the blank row 2 ("extra row") - everything below shifts up,
# hyperlinks and dimension adjust automatically.
$ws.Rows("2:2").Delete()

# Select the new row 2 (entire row) to mirror the post-delete selection state.
$ws.Range("A2:XFD2").Select()
